$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 4
$ws.Range("B2").Value = "Aircraft ActiveTrack available at max speed . When exceeding nnn, Obstacle Avoidance is not available ."
$ws.Range("C2").Value = "Aircraft ActiveTrack available at max speed"
$ws.Range("D2").Value = "0-5"
$ws.Range("E2").Value = "Missing"

$ws.Range("A3").Value = 15
$ws.Range("B3").Value = "Aircraft is close to the Home Point . Initiating Return to Home will now trigger Auto Landing ."
$ws.Range("C3").Value = "Initiating Return to Home will now trigger Auto Landing"
$ws.Range("D3").Value = "8-16"
$ws.Range("E3").Value = "Missing"

$ws.Range("A4").Value = 19
$ws.Range("B4").Value = "Aircraft is tilted , please keep the aircraft stationary and level before flight ."
$ws.Range("C4").Value = "Aircraft is tilted"
$ws.Range("D4").Value = "0-2"
$ws.Range("E4").Value = "Missing"

$ws.Range("A5").Value = 40
$ws.Range("B5").Value = "Camera error . AI Spot-Check failed . Restart camera ."
$ws.Range("C5").Value = "AI Spot-Check failed"
$ws.Range("D5").Value = "3-5"
$ws.Range("E5").Value = "Missing"

$ws.Range("A6").Value = 42
$ws.Range("B6").Value = "Camera sensor error . Hardware malfunction : Contact DJI Support to arrange for repairs ."
$ws.Range("C6").Value = "Hardware malfunction"
$ws.Range("D6").Value = "4-5"
$ws.Range("E6").Value = "Missing"

$ws.Range("A7").Value = 42
$ws.Range("B7").Value = "Camera sensor error . Hardware malfunction : Contact DJI Support to arrange for repairs ."
$ws.Range("C7").Value = "Contact DJI Support to arrange for repairs"
$ws.Range("D7").Value = "7-13"
$ws.Range("E7").Value = "Missing"

$ws.Range("A8").Value = 44
$ws.Range("B8").Value = "Cancel Return-to-Home Failed ."
$ws.Range("C8").Value = "Cancel Return-to-Home Failed"
$ws.Range("D8").Value = "0-2"
$ws.Range("E8").Value = "Missing"

$ws.Range("A9").Value = 50
$ws.Range("B9").Value = "Compass abnormal . Solution: 1. Ensure there are no magnets or metal objects near the aircraft . The ground or walls may contain metal . Move away from sources of interference before attempting flight . 2. Calibrate Compass Before Takeoff ."
$ws.Range("C9").Value = "2. Calibrate Compass Before Takeoff"
$ws.Range("D9").Value = "35-39"
$ws.Range("E9").Value = "Missing"

$ws.Range("A10").Value = 50
$ws.Range("B10").Value = "Compass abnormal . Solution: 1. Ensure there are no magnets or metal objects near the aircraft . The ground or walls may contain metal . Move away from sources of interference before attempting flight . 2. Calibrate Compass Before Takeoff ."
$ws.Range("C10").Value = "Calibrate Compass Before Takeoff"
$ws.Range("D10").Value = "36-39"
$ws.Range("E10").Value = "'False"

$ws.Range("A11").Value = 66
$ws.Range("B11").Value = "Downlink data connection lost for nnn seconds ."
$ws.Range("C11").Value = "Downlink data connection lost for nnn seconds"
$ws.Range("D11").Value = "0-6"
$ws.Range("E11").Value = "Missing"

$ws.Range("A12").Value = 66
$ws.Range("B12").Value = "Downlink data connection lost for nnn seconds ."
$ws.Range("C12").Value = "Downlink data connection lost for nnn"
$ws.Range("D12").Value = "0-5"
$ws.Range("E12").Value = "'False"

$ws.Range("A13").Value = 77
$ws.Range("B13").Value = "Exiting GPS mode : Unknown Error ."
$ws.Range("C13").Value = "Unknown Error"
$ws.Range("D13").Value = "4-5"
$ws.Range("E13").Value = "Missing"

$ws.Range("A14").Value = 81
$ws.Range("B14").Value = "Extra payload detected . Return aircraft to an area nearby the home point promptly and fly in a wind-free environment to ensure flight safety ."
$ws.Range("C14").Value = "Return aircraft to an area nearby the home point promptly and fly in a wind-free environment to ensure flight safety"
$ws.Range("D14").Value = "4-23"
$ws.Range("E14").Value = "Missing"

$ws.Range("A15").Value = 85
$ws.Range("B15").Value = "Flight altitude exceeds nnn . Aircraft may be in violation of local laws and regulations . Check and make sure you have obtained proper authorization to fly in this airspace ."
$ws.Range("C15").Value = "Aircraft may be in violation of local laws and regulations"
$ws.Range("D15").Value = "5-14"
$ws.Range("E15").Value = "Missing"

$ws.Range("A16").Value = 85
$ws.Range("B16").Value = "Flight altitude exceeds nnn . Aircraft may be in violation of local laws and regulations . Check and make sure you have obtained proper authorization to fly in this airspace ."
$ws.Range("C16").Value = "Aircraft may be in violation of local laws"
$ws.Range("D16").Value = "5-12"
$ws.Range("E16").Value = "'False"

$ws.Range("A17").Value = 91
$ws.Range("B17").Value = "GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn ."
$ws.Range("C17").Value = "GEO Zone Info: The target area is in an Altitude Zone"
$ws.Range("D17").Value = "0-10"
$ws.Range("E17").Value = "Missing"

$ws.Range("A18").Value = 91
$ws.Range("B18").Value = "GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn ."
$ws.Range("C18").Value = "GEO Zone Info:"
$ws.Range("D18").Value = "0-2"
$ws.Range("E18").Value = "'False"

$ws.Range("A19").Value = 91
$ws.Range("B19").Value = "GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn ."
$ws.Range("C19").Value = "The target area is in an Altitude Zone"
$ws.Range("D19").Value = "3-10"
$ws.Range("E19").Value = "'False"

$ws.Range("A20").Value = 92
$ws.Range("B20").Value = "GEO: You are in a Warning Zone (Airport Class Airspace Unpaved Airports Power Plant) . Fly with caution ."
$ws.Range("C20").Value = "GEO: You are in a Warning Zone (Airport Class Airspace Unpaved Airports Power Plant)"
$ws.Range("D20").Value = "0-13"
$ws.Range("E20").Value = "Missing"

$ws.Range("A21").Value = 104
$ws.Range("B21").Value = "Home Point Recorded , Return-to-Home Altitude : 98FT ."
$ws.Range("C21").Value = "Return-to-Home Altitude : 98FT"
$ws.Range("D21").Value = "4-7"
$ws.Range("E21").Value = "Missing"

$ws.Range("A22").Value = 122
$ws.Range("B22").Value = "Max Altitude Approached . Wait for the GPS satellite signal recovery before ascend ."
$ws.Range("C22").Value = "Wait for the GPS satellite signal recovery before ascend"
$ws.Range("D22").Value = "4-12"
$ws.Range("E22").Value = "Missing"

$ws.Range("A23").Value = 142
$ws.Range("B23").Value = "RTH Altitude : 98FT . Data Recorder File Index is 1 ."
$ws.Range("C23").Value = "Data Recorder File Index is 1"
$ws.Range("D23").Value = "5-10"
$ws.Range("E23").Value = "Missing"

$ws.Range("A24").Value = 146
$ws.Range("B24").Value = "SD card speed low . Change card ."
$ws.Range("C24").Value = "Change card"
$ws.Range("D24").Value = "5-6"
$ws.Range("E24").Value = "Missing"

$ws.Range("A25").Value = 159
$ws.Range("B25").Value = "Unknown Error , Cannot Takeoff . Contact DJI support ."
$ws.Range("C25").Value = "Cannot Takeoff"
$ws.Range("D25").Value = "3-4"
$ws.Range("E25").Value = "Missing"

